$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value parses as a plain number need to be forced back to
# text (matching the source inlineStr cells) since Excel auto-detects numeric
# strings typed into .Value and stores them as numbers otherwise.
$textCells = @('D5', 'D6', 'D7', 'D9', 'D10', 'D11', 'D12', 'D14', 'D16', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D27', 'D28', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D40', 'D41', 'D42', 'D43', 'D45', 'D48', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '43.791.20'
$ws.Range('E2').Value = '  -0.43%  '
$ws.Range('D3').Value = '2.288.45'
$ws.Range('E3').Value = '  -0.63%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '115.85'
$ws.Range('E5').Value = '  +1.09%  '
$ws.Range('D6').Value = '267.02'
$ws.Range('E6').Value = '  -1.38%  '
$ws.Range('D7').Value = '0.645'
$ws.Range('E7').Value = '  +2.79%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').Value = '0.617'
$ws.Range('E9').Value = '  -1.03%  '
$ws.Range('D10').Value = '48.01'
$ws.Range('E10').Value = '  -0.63%  '
$ws.Range('D11').Value = '0.0941'
$ws.Range('E11').Value = '  -1.03%  '
$ws.Range('D12').Value = '9.19'
$ws.Range('E12').Value = '  +1.00%  '
$ws.Range('E13').Value = '  +1.26%  '
$ws.Range('D14').Value = '15.53'
$ws.Range('E14').Value = '  -2.32%  '
$ws.Range('D15').Value = '2.628.80'
$ws.Range('E15').Value = '  -0.76%  '
$ws.Range('D16').Value = '0.882'
$ws.Range('E16').Value = '  +2.95%  '
$ws.Range('D17').Value = '2.282.32'
$ws.Range('E17').Value = '  -0.72%  '
$ws.Range('D18').Value = '43.649.20'
$ws.Range('E18').Value = '  -0.50%  '
$ws.Range('E19').Value = '  -0.11%  '
$ws.Range('D20').Value = '6.89'
$ws.Range('E20').Value = '  +0.96%  '
$ws.Range('D21').Value = '72.50'
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('D22').Value = '2.46'
$ws.Range('E22').Value = '  +0.24%  '
$ws.Range('D23').Value = '235.76'
$ws.Range('E23').Value = '  +0.96%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').Value = '9.58'
$ws.Range('E24').Value = '  -1.53%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').Value = '2.90'
$ws.Range('E25').Value = '  +0.10%  '
$ws.Range('E26').Value = '  +1.72%  '
$ws.Range('D27').Value = '11.73'
$ws.Range('E27').Value = '  -0.01%  '
$ws.Range('D28').Value = '42.28'
$ws.Range('E28').Value = '  +0.38%  '
$ws.Range('E29').Value = '  +0.58%  '
$ws.Range('E30').Value = '  -0.88%  '
$ws.Range('D31').Value = '174.08'
$ws.Range('E31').Value = '  -0.89%  '
$ws.Range('D32').Value = '21.77'
$ws.Range('E32').Value = '  +0.52%  '
$ws.Range('D33').Value = '0.0913'
$ws.Range('E33').Value = '  -3.12%  '
$ws.Range('D34').Value = '5.73'
$ws.Range('E34').Value = '  -0.33%  '
$ws.Range('D35').Value = '0.131'
$ws.Range('E35').Value = '  +2.04%  '
$ws.Range('D36').Value = '0.0385'
$ws.Range('E36').Value = '  +5.46%  '
$ws.Range('D37').Value = '4.69'
$ws.Range('E37').Value = '  +0.19%  '
$ws.Range('D38').Value = '3.96'
$ws.Range('E38').Value = '  +3.45%  '
$ws.Range('E39').Value = '  -1.10%  '
$ws.Range('D40').Value = '2.58'
$ws.Range('E40').Value = '  +7.55%  '
$ws.Range('D41').Value = '14.21'
$ws.Range('E41').Value = '  +2.79%  '
$ws.Range('D42').Value = '74.41'
$ws.Range('E42').Value = '  +0.46%  '
$ws.Range('D43').Value = '0.237'
$ws.Range('E43').Value = '  -2.79%  '
$ws.Range('E44').Value = '  -0.16%  '
$ws.Range('D45').Value = '5.97'
$ws.Range('E45').Value = '  -7.62%  '
$ws.Range('E46').Value = '  -2.25%  '
$ws.Range('E47').Value = '  +3.67%  '
$ws.Range('D48').Value = '8.62'
$ws.Range('E48').Value = '  -2.47%  '
$ws.Range('E49').Value = '  +0.32%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').Value = '101.59'
$ws.Range('E50').Value = '  -1.22%  '
$ws.Range('B51').Value = 'ordi'
$ws.Range('C51').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D51').Value = '72.05'
$ws.Range('E51').Value = '  +32.61%  '

# Restore the default (unstyled) cell style now that the text value is
# committed, so the written cells keep the workbook's original "no style"
# (style index 0) formatting instead of picking up a new text-format xf.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
